# Applies the "chore: update Sheets via scheduled runner" price/profit recalculation
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 1857.7778
$ws.Range("I43").Value = 2471.4285
$ws.Range("J43").Value = 1643
$ws.Range("K43").Value = 2471.4285
$ws.Range("L43").Value = 1643
$ws.Range("M43").Value = -2402.4285
$ws.Range("N43").Value = -1781

# Row 125
$ws.Range("H125").Value = 166666910
$ws.Range("I125").Value = 375
$ws.Range("J125").Value = 500000000
$ws.Range("K125").Value = 3375
$ws.Range("L125").Value = 4500000000
$ws.Range("M125").Value = -915
$ws.Range("N125").Value = -4500004920

# Row 132
$ws.Range("H132").Value = 2762.411
$ws.Range("I132").Value = 1227.9038
$ws.Range("K132").Value = 3683.7114
$ws.Range("M132").Value = -1153.7114

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 817.2
$ws.Range("I32").Value = 793.0612
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 793.0612
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -506.0612
$ws.Range("N32").Value = -2574

# Row 45
$ws.Range("H45").Value = 1271.4762
$ws.Range("I45").Value = 1099
$ws.Range("J45").Value = 1551.75
$ws.Range("K45").Value = 1099
$ws.Range("L45").Value = 1551.75
$ws.Range("M45").Value = -722
$ws.Range("N45").Value = -2305.75

# Row 107
$ws.Range("H107").Value = 20246.4
$ws.Range("J107").Value = 20246.4
$ws.Range("L107").Value = 20246.4
$ws.Range("N107").Value = -27926.4

# Row 109
$ws.Range("H109").Value = 31705
$ws.Range("J109").Value = 31705
$ws.Range("L109").Value = 31705
$ws.Range("N109").Value = -34479

# Row 132
$ws.Range("H132").Value = 2789.456
$ws.Range("I132").Value = 2640.375
$ws.Range("J132").Value = 3140.2354
$ws.Range("K132").Value = 7921.125
$ws.Range("L132").Value = 9420.706200000001
$ws.Range("M132").Value = -5391.125
$ws.Range("N132").Value = -14480.7062

$ws = $wb.Worksheets.Item("BSM")
# Row 75
$ws.Range("H75").Value = 8456.429
$ws.Range("I75").Value = 4865.8335
$ws.Range("J75").Value = 30000
$ws.Range("K75").Value = 4865.8335
$ws.Range("L75").Value = 30000
$ws.Range("M75").Value = -3929.8335
$ws.Range("N75").Value = -31872

# Row 78
$ws.Range("H78").Value = 8456.429
$ws.Range("I78").Value = 4865.8335
$ws.Range("J78").Value = 30000
$ws.Range("K78").Value = 14597.5005
$ws.Range("L78").Value = 90000
$ws.Range("M78").Value = -9917.500499999998
$ws.Range("N78").Value = -99360

# Row 134
$ws.Range("H134").Value = 2298.7793
$ws.Range("I134").Value = 2075.532
$ws.Range("J134").Value = 2798.4285
$ws.Range("K134").Value = 6226.596
$ws.Range("L134").Value = 8395.2855
$ws.Range("M134").Value = -3691.596
$ws.Range("N134").Value = -13465.2855

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2663.4033
$ws.Range("I58").Value = 2801.3674
$ws.Range("K58").Value = 2801.3674
$ws.Range("M58").Value = -2598.3674

# Row 99
$ws.Range("H99").Value = 68625.53
$ws.Range("I99").Value = 78514.08
$ws.Range("J99").Value = 4350
$ws.Range("K99").Value = 78514.08
$ws.Range("L99").Value = 4350
$ws.Range("M99").Value = -77016.08
$ws.Range("N99").Value = -7346

# Row 105
$ws.Range("H105").Value = 1235.2727
$ws.Range("I105").Value = 550
$ws.Range("J105").Value = 1387.5555
$ws.Range("K105").Value = 550
$ws.Range("L105").Value = 1387.5555
$ws.Range("M105").Value = 1197
$ws.Range("N105").Value = -4881.5555

# Row 126
$ws.Range("H126").Value = 68625.53
$ws.Range("I126").Value = 78514.08
$ws.Range("J126").Value = 4350
$ws.Range("K126").Value = 235542.24
$ws.Range("L126").Value = 13050
$ws.Range("M126").Value = -233072.24
$ws.Range("N126").Value = -17990

# Row 132
$ws.Range("H132").Value = 2403.4783
$ws.Range("I132").Value = 1285.3334
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 3856.0002
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -1326.0002
$ws.Range("N132").Value = -18560

# Row 134
$ws.Range("H134").Value = 1428.3276
$ws.Range("I134").Value = 944.5
$ws.Range("J134").Value = 2347.6
$ws.Range("K134").Value = 2833.5
$ws.Range("L134").Value = 7042.799999999999
$ws.Range("M134").Value = -298.5
$ws.Range("N134").Value = -12112.8

# Row 136
$ws.Range("H136").Value = 2663.4033
$ws.Range("I136").Value = 2801.3674
$ws.Range("K136").Value = 8404.102200000001
$ws.Range("M136").Value = -5854.102200000001

$ws = $wb.Worksheets.Item("CUL")
# Row 100
$ws.Range("H100").Value = 5602.1665
$ws.Range("J100").Value = 5817.6
$ws.Range("L100").Value = 17452.8
$ws.Range("N100").Value = -19074.8

# Row 102
$ws.Range("H102").Value = 9142.857
$ws.Range("J102").Value = 9500
$ws.Range("L102").Value = 28500
$ws.Range("N102").Value = -33368

# Row 104
$ws.Range("H104").Value = 2280.6
$ws.Range("I104").Value = 2101.5
$ws.Range("J104").Value = 2400
$ws.Range("K104").Value = 6304.5
$ws.Range("L104").Value = 7200
$ws.Range("M104").Value = -3683.5
$ws.Range("N104").Value = -12442

# Row 113
$ws.Range("H113").Value = 623
$ws.Range("I113").Value = 724
$ws.Range("J113").Value = 582.6
$ws.Range("K113").Value = 2172
$ws.Range("L113").Value = 1747.8
$ws.Range("M113").Value = -2
$ws.Range("N113").Value = -6087.8

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1430.3158
$ws.Range("I113").Value = 1041.7333
$ws.Range("J113").Value = 2887.5
$ws.Range("K113").Value = 1041.7333
$ws.Range("L113").Value = 2887.5
$ws.Range("M113").Value = 1128.2667
$ws.Range("N113").Value = -7227.5

# Row 122
$ws.Range("H122").Value = 1614.75
$ws.Range("I122").Value = 1272.125
$ws.Range("J122").Value = 2300
$ws.Range("K122").Value = 3816.375
$ws.Range("L122").Value = 6900
$ws.Range("M122").Value = -1366.375
$ws.Range("N122").Value = -11800

# Row 126
$ws.Range("H126").Value = 2171.2632
$ws.Range("I126").Value = 1410
$ws.Range("J126").Value = 3218
$ws.Range("K126").Value = 4230
$ws.Range("L126").Value = 9654
$ws.Range("M126").Value = -1760
$ws.Range("N126").Value = -14594

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2805.7693
$ws.Range("I7").Value = 2839
$ws.Range("J7").Value = 2785
$ws.Range("K7").Value = 2839
$ws.Range("L7").Value = 2785
$ws.Range("M7").Value = -2727
$ws.Range("N7").Value = -3009

# Row 40
$ws.Range("H40").Value = 3225.0908
$ws.Range("I40").Value = 2934.625
$ws.Range("J40").Value = 3999.6667
$ws.Range("K40").Value = 2934.625
$ws.Range("L40").Value = 3999.6667
$ws.Range("M40").Value = -2798.625
$ws.Range("N40").Value = -4271.6667

# Row 46
$ws.Range("H46").Value = 910.5
$ws.Range("I46").Value = 954.1
$ws.Range("J46").Value = 801.5
$ws.Range("K46").Value = 954.1
$ws.Range("L46").Value = 801.5
$ws.Range("M46").Value = -766.1
$ws.Range("N46").Value = -1177.5

# Row 61
$ws.Range("H61").Value = 1300.5
$ws.Range("I61").Value = 1226
$ws.Range("J61").Value = 1375
$ws.Range("K61").Value = 1226
$ws.Range("L61").Value = 1375
$ws.Range("M61").Value = -1024
$ws.Range("N61").Value = -1779

# Row 68
$ws.Range("H68").Value = 2879.7273
$ws.Range("I68").Value = 2320.3333
$ws.Range("J68").Value = 3089.5
$ws.Range("K68").Value = 2320.3333
$ws.Range("L68").Value = 3089.5
$ws.Range("M68").Value = -1571.3333
$ws.Range("N68").Value = -4587.5

# Row 71
$ws.Range("H71").Value = 2879.7273
$ws.Range("I71").Value = 2320.3333
$ws.Range("J71").Value = 3089.5
$ws.Range("K71").Value = 11601.6665
$ws.Range("L71").Value = 15447.5
$ws.Range("M71").Value = -7857.666499999999
$ws.Range("N71").Value = -22935.5

# Row 113
$ws.Range("H113").Value = 1300.5
$ws.Range("I113").Value = 1226
$ws.Range("J113").Value = 1375
$ws.Range("K113").Value = 1226
$ws.Range("L113").Value = 1375
$ws.Range("M113").Value = 944
$ws.Range("N113").Value = -5715

# Row 126
$ws.Range("H126").Value = 2805.7693
$ws.Range("I126").Value = 2839
$ws.Range("J126").Value = 2785
$ws.Range("K126").Value = 8517
$ws.Range("L126").Value = 8355
$ws.Range("M126").Value = -6047
$ws.Range("N126").Value = -13295

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4020.724
$ws.Range("I62").Value = 3564.3333
$ws.Range("J62").Value = 4342.8823
$ws.Range("K62").Value = 3564.3333
$ws.Range("L62").Value = 4342.8823
$ws.Range("M62").Value = -2940.3333
$ws.Range("N62").Value = -5590.8823

# Row 65
$ws.Range("H65").Value = 4020.724
$ws.Range("I65").Value = 3564.3333
$ws.Range("J65").Value = 4342.8823
$ws.Range("K65").Value = 17821.6665
$ws.Range("L65").Value = 21714.4115
$ws.Range("M65").Value = -14701.6665
$ws.Range("N65").Value = -27954.4115

# Row 107
$ws.Range("H107").Value = 243
$ws.Range("I107").Value = 237.5
$ws.Range("J107").Value = 252.16667
$ws.Range("K107").Value = 712.5
$ws.Range("L107").Value = 756.50001
$ws.Range("M107").Value = 1207.5
$ws.Range("N107").Value = -4596.50001

# Row 126
$ws.Range("H126").Value = 1650.4584
$ws.Range("I126").Value = 813.8333
$ws.Range("J126").Value = 2487.0833
$ws.Range("K126").Value = 2441.4999
$ws.Range("L126").Value = 7461.249899999999
$ws.Range("M126").Value = 28.5001000000002
$ws.Range("N126").Value = -12401.2499
